# Add 7 new data rows (457-463) to the Landscaping Data sheet, continuing
# the existing table that ended at row 456.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Copy the date-format style from the last existing row (A456) down onto
#    the new A457:A463 cells before writing values, so the new cells keep
#    the same numFmt-14 ("m/d/yyyy") style used throughout column A.
# ---------------------------------------------------------------------------
$null = $ws.Range("A456").Copy()
$null = $ws.Range("A457:A463").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# 2. Row data for the 7 new records (columns A-E, G-T; F is a formula).
#    Array order matches the sheet columns, skipping F (Temp_Diff):
#      0:A Date, 1:B Plant_Type, 2:C Plant_Size, 3:D Low, 4:E High,
#      5:G Rain, 6:H Growth, 7:I Pruned, 8:J Quadrant, 9:K Shade, 10:L UV,
#      11:M Humidity, 12:N Dew_Point, 13:O Pressure, 14:P Wind_Gust,
#      15:Q Cloud_Cover, 16:R Visibility, 17:S AQI, 18:T Pollen
# ---------------------------------------------------------------------------
$rows = @{
    457 = @(45852, "Flowering",     "Large",  69, 82, 0,    0.1, "No", 2, "Bright",  8, 0.69, 71, 30.04, 5, 0.64, 8.7, 41, 0)
    458 = @(45852, "Nonflowering",  "Medium", 69, 82, 0,    0.1, "No", 3, "Bright",  8, 0.69, 71, 30.04, 5, 0.64, 8.7, 41, 0)
    459 = @(45852, "Nonflowering",  "Small",  69, 82, 0,    0.2, "No", 3, "Bright",  8, 0.69, 71, 30.04, 5, 0.64, 8.7, 41, 0)
    460 = @(45852, "Nonflowering",  "Medium", 69, 82, 0,    0,   "No", 3, "Neutral", 8, 0.69, 71, 30.04, 5, 0.64, 8.7, 41, 0)
    461 = @(45852, "Nonflowering",  "Medium", 69, 82, 0,    0,   "No", 3, "Neutral", 8, 0.69, 71, 30.04, 5, 0.64, 8.7, 41, 0)
    462 = @(45852, "Nonflowering",  "Large",  69, 82, 0,    0.3, "No", 4, "Bright",  8, 0.69, 71, 30.04, 5, 0.64, 8.7, 41, 0)
    463 = @(45852, "Tree",          "Medium", 69, 82, 0,    0.5, "No", 1, "Neutral", 8, 0.69, 71, 30.04, 5, 0.64, 8.7, 41, 0)
}

foreach ($r in 457..463) {
    $vals = $rows[$r]

    $ws.Cells.Item($r, 1).Value  = $vals[0]   # A Date
    $ws.Cells.Item($r, 2).Value  = $vals[1]   # B Plant_Type
    $ws.Cells.Item($r, 3).Value  = $vals[2]   # C Plant_Size
    $ws.Cells.Item($r, 4).Value  = $vals[3]   # D Low
    $ws.Cells.Item($r, 5).Value  = $vals[4]   # E High
    # F Temp_Diff is filled in below as a shared formula
    $ws.Cells.Item($r, 7).Value  = $vals[5]   # G Rain
    $ws.Cells.Item($r, 8).Value  = $vals[6]   # H Growth
    $ws.Cells.Item($r, 9).Value  = $vals[7]   # I Pruned
    $ws.Cells.Item($r, 10).Value = $vals[8]   # J Quadrant
    $ws.Cells.Item($r, 11).Value = $vals[9]   # K Shade
    $ws.Cells.Item($r, 12).Value = $vals[10]  # L UV
    $ws.Cells.Item($r, 13).Value = $vals[11]  # M Humidity
    $ws.Cells.Item($r, 14).Value = $vals[12]  # N Dew_Point
    $ws.Cells.Item($r, 15).Value = $vals[13]  # O Pressure
    $ws.Cells.Item($r, 16).Value = $vals[14]  # P Wind_Gust
    $ws.Cells.Item($r, 17).Value = $vals[15]  # Q Cloud_Cover
    $ws.Cells.Item($r, 18).Value = $vals[16]  # R Visibility
    $ws.Cells.Item($r, 19).Value = $vals[17]  # S AQI
    $ws.Cells.Item($r, 20).Value = $vals[18]  # T Pollen
}

# ---------------------------------------------------------------------------
# 3. Fill column F (Temp_Diff = ABS(D-E)) for the new rows with one formula
#    assignment so it is written as a single shared formula group.
# ---------------------------------------------------------------------------
$ws.Range("F457:F463").Formula = "=ABS(D457-E457)"

# ---------------------------------------------------------------------------
# 4. Leave the selection where the author left it after entering the data.
# ---------------------------------------------------------------------------
$null = $ws.Range("I464").Select()

Write-Output "Added rows 457-463"
